# Apply LaTeX-ified updates to the "key_buildingBlock_pairs" workbook.
# Sheet1 holds (key, description, building-block-label) triples; this edit
# wraps several numeric/variable tokens in LaTeX math delimiters ($...$)
# and adds a couple of brand-new building-block label cells in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New cells (column C was empty for these rows before) -----------------
$ws.Range("C1").Value = "`$a_{10}`$; `$k`$;"
$ws.Range("C2").Value = "정수 `$-7<  k<  20`$의 개수; "
$ws.Range("C20").Value = "`$\sin \theta+\cos \theta`$;"
$ws.Range("C45").Value = "`$\displaystyle\lim _{x \rightarrow-1-} f(x)+\displaystyle\lim _{x \rightarrow 2} f(x)`$"

# --- Existing cells whose text changed -------------------------------------
$ws.Range("B16").Value = "주어진 삼각방정식을 좌변 인수분해, 우변 `$0`$이 되도록 변형합니다."
$ws.Range("C17").Value = "`$\tan \theta`$;"
$ws.Range("C18").Value = "`$\tan \theta`$에 대한 2차식 인수분해;"
$ws.Range("B32").Value = "구한 두 점의 `$y`$ 좌표가 일치하도록 방정식을 세웁니다."
$ws.Range("C42").Value = "`$f^{\prime}(1)`$; "
$ws.Range("C46").Value = "좌변 삼차함수, 우변 `$k`$;"
$ws.Range("C52").Value = "`$0`$에서 `$k`$ 까지;"
$ws.Range("B53").Value = "우변`$0`$인 항등식의 좌변을 인수분해합니다."
$ws.Range("C55").Value = "최대 `$1`$ 최소 `$0`$;"

# --- Selection / window cosmetics ------------------------------------------
$ws.Range("B10").Select()
